$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Valor Mora" total (E11): 90189 -> 176820
$ws.Range("E11").Value = 176820

# 2. Update existing data row 16 (was worker 20177889 / PAULA MELISSA MOUTHON GARCIA / 1812)
#    to worker 1001970789 / DANNA MILENA DIAZ SOLANO / period 2503, with refreshed amounts
$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "1001970789"
$ws.Cells.Item(16, 4).Value = "DANNA MILENA DIAZ SOLANO"
$ws.Cells.Item(16, 5).Value = "2503"
$ws.Cells.Item(16, 6).Value = 58940
$ws.Cells.Item(16, 7).Value = 1473500

# 3. Insert a brand-new row at 17 (pushes the old row 17 and the footer rows down by one)
$ws.Rows.Item(17).Insert()

# 4. Clone the formatting of row 16 onto the freshly inserted row 17
$src = $ws.Range("B16:J16")
$dst = $ws.Range("B17:J17")
$src.Copy()
$dst.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# 5. Populate new row 17: same worker, new period 2502
$ws.Cells.Item(17, 2).Value = "CC"
$ws.Cells.Item(17, 3).Value = "1001970789"
$ws.Cells.Item(17, 4).Value = "DANNA MILENA DIAZ SOLANO"
$ws.Cells.Item(17, 5).Value = "2502"
$ws.Cells.Item(17, 6).Value = 58940
$ws.Cells.Item(17, 7).Value = 1473500

# 6. Row 18 (the old row 17, shifted down by the insert) becomes a new worker entry
$ws.Cells.Item(18, 2).Value = "CC"
$ws.Cells.Item(18, 3).Value = "1047426460"
$ws.Cells.Item(18, 4).Value = "YINA DANNELYS MARTINEZ BARRIOS"
$ws.Cells.Item(18, 5).Value = "2502"
$ws.Cells.Item(18, 6).Value = 58940
$ws.Cells.Item(18, 7).Value = 1473500
